$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Table 1 (columns A-D): Los Angeles Dodgers righties ---
$ws1.Range("A1").Value = "Los Angeles Dodgers righties"

$ws1.Range("A3").Value = "Taylor"
$ws1.Range("B3").Value = 3300
$ws1.Range("C3").Value = 19

$ws1.Range("A4").Value = "Verdugo"
$ws1.Range("B4").Value = 2500
$ws1.Range("C4").Value = 0

$ws1.Range("A5").Value = "Kemp"
$ws1.Range("B5").Value = 2800
$ws1.Range("C5").Value = 28.2

$ws1.Range("A6").Value = "Pederson"
$ws1.Range("B6").Value = 2200
$ws1.Range("C6").Value = 12.2

# --- Table 2 (columns F-I): Minnesota Twins hitters ---
$ws1.Range("F1").Value = "Minnesota Twins hitters"

$ws1.Range("F3").Value = "Mauer"
$ws1.Range("G3").Value = 3200
$ws1.Range("H3").Value = 18.4

$ws1.Range("F4").Value = "Dozier"
$ws1.Range("G4").Value = 3300
$ws1.Range("H4").Value = 3

$ws1.Range("F5").Value = "Garver"
$ws1.Range("G5").Value = 2100
$ws1.Range("H5").Value = 9.5

$ws1.Range("F6").Value = "Morrison"
$ws1.Range("G6").Value = 2600
$ws1.Range("H6").Value = 6.2

# --- Table 3 (columns K-N): Colorado Rockies righties ---
$ws1.Range("K1").Value = "Colorado Rockies righties"

$ws1.Range("K3").Value = "Desmond"
$ws1.Range("L3").Value = 2700
$ws1.Range("M3").Value = 0

$ws1.Range("K4").Value = "Story"
$ws1.Range("L4").Value = 3600
$ws1.Range("M4").Value = 9

$ws1.Range("K5").Value = "Iannetta"
$ws1.Range("L5").Value = 2300
$ws1.Range("M5").Value = 9.2

$ws1.Range("K6").Value = "Cuevas"
$ws1.Range("L6").Value = 2200
$ws1.Range("M6").Value = 3

# Table 2 "Success"/"Failure" banner flips from Success to Failure
$ws1.Range("I8").Value = "Failure"

# --- Selection changes ---
$ws1.Activate()
$ws1.Range("F9:I9").Select()

$ws2.Activate()
$ws2.Range("A1:A3").Select()

$ws1.Activate()
